{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright footer block that\n// Jekyll appended after the \"Requisitos\" section (and the blank paragraph\n// that separated it from the LOT2038 requisito line).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOT2038: ...\" requirement paragraph - the footer block we\n// need to remove always starts immediately after it.\nlet lot2038Index = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOT2038\") !== -1) {\n    lot2038Index = i;\n    break;\n  }\n}\n\nif (lot2038Index !== -1) {\n  // Delete the paragraphs that follow LOT2038 up to (and including) the\n  // copyright/footer line, stopping right before the next blank paragraph\n  // that precedes the final page-break paragraph.\n  const toDelete = [];\n  for (let i = lot2038Index + 1; i < items.length; i++) {\n    const text = items[i].text.trim();\n    toDelete.push(items[i]);\n    if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n      break;\n    }\n  }\n  // Delete from the end backwards so indices/objects stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright footer block that\n# Jekyll appended after the \"Requisitos\" section (and the blank paragraph\n# that separated it from the LOT2038 requisito line).\n$d = $word.ActiveDocument\n\n# Locate the \"LOT2038: ...\" requirement paragraph - the footer block we\n# need to remove always starts immediately after it.\n$lot2038Index = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOT2038*\") {\n        $lot2038Index = $i\n        break\n    }\n}\n\nif ($lot2038Index -ne -1) {\n    # Walk forward from the paragraph right after LOT2038 until (and\n    # including) the copyright/footer line.\n    $endIndex = -1\n    for ($i = $lot2038Index + 1; $i -le $d.Paragraphs.Count; $i++) {\n        $text = $d.Paragraphs.Item($i).Range.Text\n        if ($text -like \"*Powered by Jekyll*\") {\n            $endIndex = $i\n            break\n        }\n    }\n\n    if ($endIndex -ne -1) {\n        $startPara = $d.Paragraphs.Item($lot2038Index + 1)\n        $endPara = $d.Paragraphs.Item($endIndex)\n        $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n        $range.Delete()\n    }\n}\n"}
